$wb = $excel.ActiveWorkbook

# 1. Update the Date value on the Metadata sheet
$wsMeta = $wb.Worksheets.Item("Metadata")
$wsMeta.Range("B8").Value2 = "2024-03-19T13:17:15+00:00"

# 2. Swap columns AK (37) and AL (38) on the Elements sheet:
#    the "Mapping: RIM Mapping" column and the
#    "Mapping: Spécification métier vers l'extension ROR Comment"
#    column traded places (header text + each row's data).
$ws = $wb.Worksheets.Item("Elements")

$rows = @(1, 3, 5, 6)
foreach ($row in $rows) {
    $akCell = $ws.Range("AK$row")
    $alCell = $ws.Range("AL$row")
    $akVal = $akCell.Value2
    $alVal = $alCell.Value2
    $akCell.Value2 = $alVal
    $alCell.Value2 = $akVal
}

# 3. Swap the column widths to match the new content
#    (AK becomes the wide column, AL becomes the narrow one).
$ws.Columns.Item(37).ColumnWidth = 64.0
$ws.Columns.Item(38).ColumnWidth = 24.166666666666664
